# Update column F ("dSF") values to reflect the repulled data / mean
# calculation changes described in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 4
    5  = -1
    10 = -1
    12 = 0
    14 = 1
    16 = -8
    22 = 5
    27 = 8
    28 = 3
    29 = 5
    30 = 4
    32 = 0
    35 = 10
    38 = -4
    40 = 8
    44 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
